$d = $word.ActiveDocument

# Contact-info line to insert, centered, directly under the name/title paragraph.
$contactText = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Paragraph 3 ("Senior data scientist...") uses the plain/default paragraph style with
# no direct run formatting, so splitting right before it yields a brand-new, completely
# unformatted paragraph (no pStyle, no rPr) - exactly what the target diff wants, rather
# than one that inherits the bold/large-font run formatting of the name line above it.
$bodyPara = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($bodyPara.Range.Start, $bodyPara.Range.Start)
$insertPoint.InsertBefore($contactText + [char]13)

# The freshly created paragraph is now paragraph 3; center it (adds only <w:jc val="center"/>).
$newPara = $d.Paragraphs.Item(3)
$newPara.Range.ParagraphFormat.Alignment = 1

# Move that whole clean paragraph (text + its own paragraph mark) up so it sits right
# after the "Dheeraj Chand" title paragraph (paragraph 1), before "PROFESSIONAL SUMMARY".
$newPara.Range.Cut()
$titleEnd = $d.Paragraphs.Item(1).Range.End
$d.Range($titleEnd, $titleEnd).Paste()
